$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing dates in rows 84 & 85 ---
$ws.Range("A84").Value = 43460
$ws.Range("A85").Value = 43462

# --- Add new row 86: club statistics entry ---
# Copy the date-number-format from the row above so the new cell reuses
# the existing style (numFmtId 14 date format) instead of minting a new one.
$ws.Range("A85").Copy() | Out-Null
$ws.Range("A86").PasteSpecial(-4122) | Out-Null

$ws.Range("A86").Value = 43463
$ws.Range("B86").Value = "Daniel"
$ws.Range("C86").Value = "Code improvement, pdf Club"
$ws.Range("D86").Value = 5

# --- Update selection to reflect the newly active cell ---
$ws.Range("A84").Select() | Out-Null
